$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("D27").Formula = "=D20+D21+D22+D23-D24+D25-D26"

$ws.Activate()
$ws.Range("D4:G4").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$win.Zoom = 59
